# feat(CMM): part provider mod hints & cleanups
#
# - Renames the "ccm_ui_remove" localization id to "cmm_ui_remove"
# - Adds a new "cmm_ui_unknown" localization row (JP: 不明なMOD / CN: 未知Mod)
# - Restyles row 4 to the plain style used elsewhere in the sheet
# - Normalizes the JP (微软雅黑) / CN (宋体, no explicit theme color) hint fonts
#   that now live on the new row
# - Updates the saved selection to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clean up the JP/CN hint fonts on row 4 *before* we move anything around,
#    so that the new "unknown mod" row below inherits the cleaned-up fonts.
# ---------------------------------------------------------------------------
$ws.Range("C4").Font.Name = "微软雅黑"   # was "Yu Gothic"
$ws.Range("D4").Font.Color = 0           # drop the explicit theme color on 宋体

# ---------------------------------------------------------------------------
# 2. Flatten the (empty) B column cell of row 4 to the plain column-A style
#    first, so that when row 4 is duplicated below both B4 and B5 end up on
#    the same plain style.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Duplicate row 4 (with its now-updated fonts) down into row 5 - this is
#    the new "cmm_ui_unknown" hint row - without disturbing any later rows.
# ---------------------------------------------------------------------------
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Re-flatten row 4's C4/D4 cells back to the plain style used by column A
#    / the rows above it, now that row 5 has its own copy of the JP/CN
#    fonts.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Fill in the actual text. Row 4 becomes the renamed "cmm_ui_remove" hint,
#    row 5 becomes the new "cmm_ui_unknown" hint.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value2 = "cmm_ui_remove"
$ws.Range("C4").Value2 = "無し"
$ws.Range("D4").Value2 = "空白"

$ws.Range("A5").Value2 = "cmm_ui_unknown"
$ws.Range("C5").Value2 = "不明なMOD"
$ws.Range("D5").Value2 = "未知Mod"

# ---------------------------------------------------------------------------
# 6. Row heights: row 4 now matches the lighter rows above it, row 5 matches
#    row 4.
# ---------------------------------------------------------------------------
$ws.Rows(4).RowHeight = 23.25
$ws.Rows(5).RowHeight = 23.25

# ---------------------------------------------------------------------------
# 7. Restore the saved cursor position to C9.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select() | Out-Null
